$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the next day's row to the bottom of the "Gold data" table.
# Column A gets the new date (a brand-new shared string).
$ws.Cells.Item(95, 1).Value = "20-12-2025"

# Column B reuses the same price text as the row above (no new price was
# scraped for this date), so copy the cell instead of retyping the text -
# this makes the engine reuse the existing shared-string entry rather than
# creating a duplicate one.
$ws.Range("B94").Copy($ws.Range("B95"))
